$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "EXCVMUX" label to "EXCVectMUX"
$ws.Range("BC1").Value = " EXCVectMUX"

# Remove the "Interrupt Enable" bit columns (LD.INT_EN / INT_ENMUX header + data)
$ws.Range("BE1:BF65").ClearContents()
